$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Fix the "data-sharing-manifesto" bookmark so that it wraps only the
# heading's run text ("Data Sharing Manifesto") instead of collapsing at
# the very start of the paragraph. Concretely: bookmarkStart stays right
# before the run, and bookmarkEnd moves to right after the run (both
# still inside the Heading 1 paragraph, before the paragraph mark).
# ---------------------------------------------------------------------------
$bookmarkName = "data-sharing-manifesto"

if ($d.Bookmarks.Exists($bookmarkName)) {
    $existing = $d.Bookmarks.Item($bookmarkName)
    # Anchor on the paragraph that currently contains the bookmark, then
    # trim the trailing paragraph mark so the bookmark only covers the
    # visible heading text, matching how Word scopes a heading bookmark to
    # the run(s) rather than the whole paragraph (incl. pilcrow).
    $headingRange = $existing.Range.Paragraphs(1).Range.Duplicate
    $headingRange.MoveEnd(1, -1) | Out-Null

    $existing.Delete()
    $d.Bookmarks.Add($bookmarkName, $headingRange) | Out-Null
}
